$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" sending-cluster block (old rows 8-10) was removed entirely
# as part of the TPM recompute; deleting these rows also shifts rows 2-7 up
# (they stay put here, only 8-10 go) and drops the now-unused "MuSCs" shared string.
$ws.Rows("8:10").Delete()

# Row 2 updated values
$ws.Range("I2").Value = 0.02394963654761903
$ws.Range("J2").Value = 0.02394963654761903
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.058783666666667
$ws.Range("N2").Value = 3.176351
$ws.Range("O2").Value = 0.5014862149947701
$ws.Range("P2").Value = 0.5014862149947702
$ws.Range("Q2").Value = 0.9075343023160001
$ws.Range("R2").Value = 8.167808720844
$ws.Range("S2").Value = 0.01201041258276588
$ws.Range("T2").Value = 0.01201041258276588

# Row 3 updated values
$ws.Range("I3").Value = 0.02394963654761903
$ws.Range("J3").Value = 0.02394963654761903
$ws.Range("O3").Value = 0.475327031240749
$ws.Range("P3").Value = 0.4753270312407492
$ws.Range("S3").Value = 0.01138390963947469
$ws.Range("T3").Value = 0.0113839096394747

# Row 4 updated values
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("I4").Value = 0.02394963654761903
$ws.Range("J4").Value = 0.02394963654761903
$ws.Range("M4").Value = 0.048954
$ws.Range("N4").Value = 0.146862
$ws.Range("O4").Value = 0.02318675376448066
$ws.Range("P4").Value = 0.02318675376448067
$ws.Range("Q4").Value = 0.041960823192
$ws.Range("R4").Value = 0.377647408728
$ws.Range("S4").Value = 0.0005553143253784491
$ws.Range("T4").Value = 0.0005553143253784492

# Row 5 updated values
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 34.93245566666667
$ws.Range("H5").Value = 104.797367
$ws.Range("I5").Value = 0.9760503634523809
$ws.Range("J5").Value = 0.9760503634523809
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.058783666666667
$ws.Range("N5").Value = 3.176351
$ws.Range("O5").Value = 0.5014862149947701
$ws.Range("P5").Value = 0.5014862149947702
$ws.Range("Q5").Value = 36.98591349642412
$ws.Range("R5").Value = 332.8732214678171
$ws.Range("S5").Value = 0.4894758024120042
$ws.Range("T5").Value = 0.4894758024120043

# Row 6 updated values
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 34.93245566666667
$ws.Range("H6").Value = 104.797367
$ws.Range("I6").Value = 0.9760503634523809
$ws.Range("J6").Value = 0.9760503634523809
$ws.Range("O6").Value = 0.475327031240749
$ws.Range("P6").Value = 0.4753270312407492
$ws.Range("Q6").Value = 35.056605614106
$ws.Range("R6").Value = 315.509450526954
$ws.Range("S6").Value = 0.4639431216012743
$ws.Range("T6").Value = 0.4639431216012744

# Row 7 updated values
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 34.93245566666667
$ws.Range("H7").Value = 104.797367
$ws.Range("I7").Value = 0.9760503634523809
$ws.Range("J7").Value = 0.9760503634523809
$ws.Range("M7").Value = 0.048954
$ws.Range("N7").Value = 0.146862
$ws.Range("O7").Value = 0.02318675376448066
$ws.Range("P7").Value = 0.02318675376448067
$ws.Range("Q7").Value = 1.710083434706
$ws.Range("R7").Value = 15.390750912354
$ws.Range("S7").Value = 0.02263143943910221
$ws.Range("T7").Value = 0.02263143943910222
